# Publicação da Estrutura do Header
# In the Sprint Backlog table, the row describing the "Header" task had an
# empty "Conclusão" (last) cell. Mark it as concluded by writing "SIM" into
# that cell.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$numCols = $tbl.Columns.Count
$targetRow = $null

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $taskCell = $tbl.Cell($r, 1)
    $taskText = $taskCell.Range.Text
    if ($taskText -match "Header") {
        $targetRow = $r
        break
    }
}

if ($targetRow -ne $null) {
    $conclusionCell = $tbl.Cell($targetRow, $numCols)
    $conclusionCell.Range.Text = "SIM"
}
